# Reposition the four "Phase activities" content placeholders on slide 1
# (shape ids 17, 18, 19, 20 / placeholder idx 23, 24, 25, 26). They were
# inheriting their geometry from the layout; the edit pins each one to an
# explicit position (moved further up the slide) while keeping the same
# width/height.
#
# NOTE: PowerPoint's COM Shape.Left/Top/Width/Height are expressed in
# points and are stored internally as single-precision (32-bit) floats,
# so the literals below are chosen so that point -> EMU (1 pt = 12700 EMU)
# conversion reproduces the exact target EMU offsets.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targets = @(
    @{ Id = 17; Left = 53.75;              Top = 122.08716535433071; Width = 192.75; Height = 370.18023622047247 },
    @{ Id = 18; Left = 273.9250183105469;  Top = 122.08716535433071; Width = 192.75; Height = 370.18023622047247 },
    @{ Id = 19; Left = 493.32501220703125; Top = 121.62889763779528; Width = 192.75; Height = 370.18023622047247 },
    @{ Id = 20; Left = 711.7133070866141;  Top = 122.0809097290039;  Width = 192.75; Height = 370.18023622047247 }
)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    foreach ($t in $targets) {
        if ($sh.Id -eq $t.Id) {
            $sh.Left = $t.Left
            $sh.Top = $t.Top
            $sh.Width = $t.Width
            $sh.Height = $t.Height
        }
    }
}
